$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.969.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.78%  '
$ws.Range("D3").Value = "'1.866.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.22%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'318.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = "'0.4367"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.05%  '
$ws.Range("D8").Value = "'0.3731"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.65%  '
$ws.Range("D9").Value = "'0.07479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("D10").Value = "'0.9366"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'21.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.83%  '
$ws.Range("D12").Value = "'1.935.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = "'6.731"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").Value = "'5.442"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.46%  '
$ws.Range("D15").Value = "'0.06887"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = "'81.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.57%  '
$ws.Range("D18").Value = "'0.000009037"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.40%  '
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  -5.07%  '
$ws.Range("D21").Value = "'27.944.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.82%  '
$ws.Range("D22").Value = "'5.123"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.17%  '
$ws.Range("D23").Value = "'11.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = "'2.144.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").Value = "'2.038"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("D26").Value = "'154.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = "'5.559"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = "'113.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.87%  '
$ws.Range("D30").Value = "'1.705"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.26%  '
$ws.Range("D31").Value = "'0.09039"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").Value = "'0.8174"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.86%  '
$ws.Range("D33").Value = "'4.821"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.90%  '
$ws.Range("D34").Value = "'1.176"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.50%  '
$ws.Range("D35").Value = "'2.970"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("D36").Value = "'1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").Value = "'0.05526"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("D38").Value = "'1.123"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("D39").Value = "'0.01979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.53%  '
$ws.Range("D40").Value = "'2.950"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("D41").Value = "'0.5269"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'7.038"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.80%  '
$ws.Range("E43").Value = '  -3.14%  '
$ws.Range("D44").Value = "'8.795"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.54%  '
$ws.Range("D45").Value = "'0.06746"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.40%  '
$ws.Range("D46").Value = "'0.4897"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.91%  '
$ws.Range("E47").Value = '  -6.27%  '
$ws.Range("D48").Value = "'107.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = "'1.912"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -14.29%  '
$ws.Range("D50").Value = "'1.002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("D51").Value = "'1.675"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.09%  '
